$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / safe updates (not ambiguous with numbers) ---
$ws.Range("D2").Value = "42.888.27"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.213.75"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("E6").Value = "  +3.94%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "2.547.66"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "2.216.85"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "42.848.82"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  +4.41%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -4.68%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  +8.60%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +9.37%  "
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +7.30%  "
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("E40").Value = "  +17.68%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("B46").Value = "WOONetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("E51").Value = "  +22.50%  "

# --- Numeric-looking price text that must stay as text (e.g. "7.00", "0.615") ---
# Format each area as Text before assignment (loop over Areas: a union Range
# only applies NumberFormat to its first Area), then clear the number-format
# afterwards so style attributes match the original (unstyled) cells.
$numRange = $ws.Range("D5,D6,D7,D9,D10,D11,D12,D15,D20,D21,D22,D23,D24,D26,D27,D29,D30,D31,D32,D33,D36,D38,D39,D40,D43,D45,D46,D47,D48,D50,D51")
foreach ($area in $numRange.Areas) { $area.NumberFormat = "@" }
$ws.Range("D5").Value = "257.22"
$ws.Range("D6").Value = "77.97"
$ws.Range("D7").Value = "0.615"
$ws.Range("D9").Value = "0.591"
$ws.Range("D10").Value = "43.16"
$ws.Range("D11").Value = "0.0906"
$ws.Range("D12").Value = "7.00"
$ws.Range("D15").Value = "14.43"
$ws.Range("D20").Value = "71.07"
$ws.Range("D21").Value = "5.98"
$ws.Range("D22").Value = "2.29"
$ws.Range("D23").Value = "229.73"
$ws.Range("D24").Value = "9.29"
$ws.Range("D26").Value = "42.70"
$ws.Range("D27").Value = "10.73"
$ws.Range("D29").Value = "2.20"
$ws.Range("D30").Value = "2.21"
$ws.Range("D31").Value = "173.54"
$ws.Range("D32").Value = "20.40"
$ws.Range("D33").Value = "0.0874"
$ws.Range("D36").Value = "0.0355"
$ws.Range("D38").Value = "4.39"
$ws.Range("D39").Value = "13.19"
$ws.Range("D40").Value = "2.85"
$ws.Range("D43").Value = "61.21"
$ws.Range("D45").Value = "103.42"
$ws.Range("D46").Value = "0.483"
$ws.Range("D47").Value = "8.44"
$ws.Range("D48").Value = "0.0973"
$ws.Range("D50").Value = "1.14"
$ws.Range("D51").Value = "1.48"
foreach ($area in $numRange.Areas) { $area.ClearFormats() }
